$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.284.12"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "1.785.70"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'338.32"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.3830"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("D8").Value = "'0.3438"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'46.94"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'1.151"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").Value = "'0.07391"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "'23.41"
$ws.Range("E12").Value = "  +7.72%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'6.459"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "'7.326"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "1.783.86"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'0.00001079"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'82.04"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'6.438"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "28.251.33"
$ws.Range("E23").Value = "  +3.80%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "'2.359"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'1.437"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").Value = "'20.72"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'2.414"
$ws.Range("E28").Value = "  -3.85%  "
$ws.Range("D29").Value = "'154.56"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "1.986.61"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'134.79"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'4.014"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "'6.094"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").Value = "'0.08882"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'12.75"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'0.02410"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'0.6857"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'5.350"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "'0.06405"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'0.2166"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "'1.249"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "'1.501"
$ws.Range("E42").Value = "  -6.79%  "
$ws.Range("D43").Value = "'8.249"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Value = "'14.20"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'0.6314"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'3.876"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'133.53"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").Value = "'2.080"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").Value = "'0.07509"
$ws.Range("E50").Value = "  +5.60%  "
$ws.Range("D51").Value = "'1.210"
$ws.Range("E51").Value = "  +6.34%  "
